$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Processes sheet ("Hoja4"): rotate the fuel/product/type columns so
# that the "type" column (PRODUCTIVE/DISSIPATIVE, validated against
# Validate!$B$2:$B$3) now sits in column B, "fuel" moves to column C
# and "product" moves to column D.
# -----------------------------------------------------------------
$wsProc = $wb.Worksheets.Item("Processes")

$wsProc.Range("B1").Value = "type"
$wsProc.Range("C1").Value = "fuel"
$wsProc.Range("D1").Value = "product"

$wsProc.Range("B2").Value = "PRODUCTIVE"
$wsProc.Range("C2").Value = "B1"
$wsProc.Range("D2").Value = "B2+B3"

$wsProc.Range("B3").Value = "PRODUCTIVE"
$wsProc.Range("C3").Value = "B2"
$wsProc.Range("D3").Value = "B4"

$wsProc.Range("B4").Value = "PRODUCTIVE"
$wsProc.Range("C4").Value = "B3"
$wsProc.Range("D4").Value = "B5"

# Column B is now wider to fit "PRODUCTIVE"/"DISSIPATIVE".
$wsProc.Columns.Item(2).ColumnWidth = 11.5

# The PRODUCTIVE/DISSIPATIVE list validation follows the "type" column
# from D to B.
$wsProc.Range("D2:D4").Validation.Delete()
$wsProc.Range("B2:B4").Validation.Add(3, 1, 1, "=Validate!B2:B3")

# -----------------------------------------------------------------
# View/selection state: ResourcesCost used to be the active tab with
# D5 selected; now Processes is active (with B1:B4 selected) and
# ResourcesCost is left with C2 selected.
# -----------------------------------------------------------------
$wsCost = $wb.Worksheets.Item("ResourcesCost")
$wsCost.Activate()
$wsCost.Range("C2").Select()

$wsProc.Activate()
$wsProc.Range("B1:B4").Select()
